$wb = $excel.ActiveWorkbook
$ch = $wb.Worksheets.Item("Chests")
$ch.Columns.Item(6).Width = 22
$ch.Columns.Item(7).Width = 98.140625
Write-Host "done"
